$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark "Accepted" in the Acceptance column (D) for all rows except row 7 (D7 stays empty)
$ws.Range("D2").Value = "Accepted"
$ws.Range("D3").Value = "Accepted"
$ws.Range("D4").Value = "Accepted"
$ws.Range("D5").Value = "Accepted"
$ws.Range("D6").Value = "Accepted"
$ws.Range("D8").Value = "Accepted"
$ws.Range("D9").Value = "Accepted"
$ws.Range("D10").Value = "Accepted"
$ws.Range("D11").Value = "Accepted"

# Row 5: add a follow-up comment in column F (wrapped, left/top aligned)
$ws.Range("F5").Value = "FYI: it shall be TiRight_u8GetStatus but it's a `nminor point"
$ws.Range("F5").WrapText = $true
$ws.Range("F5").HorizontalAlignment = -4131
$ws.Range("F5").VerticalAlignment = -4160

# Row 7: point status changes to Open, with a new comment
$ws.Range("E7").Value = "Open"
$ws.Range("F7").Value = "for u8Pin range shall be --> 0 - 43 (put here the max num of pins)`nAlso LED_1, LED_2 point still open"

# Reuse F5's newly-created style for F7 via copy/paste-special (format only) so the
# stylesheet doesn't accumulate transient/unused cellXf combinations
$ws.Range("F5").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to match final state
$ws.Range("D9").Select()
